$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated counts ("Relatório com dados de 17/12/2021")
$ws.Range("B104").Value = 7975
$ws.Range("B127").Value = 3858
$ws.Range("B250").Value = 420184
$ws.Range("B251").Value = 261376
$ws.Range("B269").Value = 79555
$ws.Range("B270").Value = 160742
$ws.Range("B302").Value = 219
$ws.Range("B335").Value = 50
$ws.Range("B360").Value = 187
$ws.Range("B362").Value = 5636
$ws.Range("B479").Value = 6654
$ws.Range("B553").Value = 7
$ws.Range("B654").Value = 1248
$ws.Range("B698").Value = 3224
$ws.Range("B709").Value = 290
$ws.Range("B735").Value = 430
$ws.Range("B753").Value = 208895
$ws.Range("B758").Value = 936
$ws.Range("B871").Value = 4042
$ws.Range("B938").Value = 14563
$ws.Range("B946").Value = 9755
$ws.Range("B957").Value = 5558
$ws.Range("B979").Value = 62228
$ws.Range("B1001").Value = 13
$ws.Range("B1036").Value = 3
$ws.Range("B1043").Value = 50330
$ws.Range("B1049").Value = 1003
$ws.Range("B1050").Value = 380
$ws.Range("B1065").Value = 105
$ws.Range("B1075").Value = 16
$ws.Range("B1085").Value = 47

# New row 1087: a new "tema" (STF RG 1189) with 1 stayed process
$ws.Range("A1087").Value = "STF RG 1189"
$ws.Range("B1087").Value = 1

# Copy formatting (fill/border/font/number format) from the last existing
# data row so the new row matches the table's established style.
$ws.Range("A1086:B1086").Copy()
$ws.Range("A1087:B1087").PasteSpecial(-4122)
$ws.Rows.Item(1087).RowHeight = $ws.Rows.Item(1086).RowHeight
